$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fila 25: el estado de la tarea "sistematizar el RG 14..." pasa de
# "no comenzado" a "en proceso"
$ws.Range("B25").Value = "en proceso"

# Nueva tarea (fila 26): formulario de bancos - baja logica de cuentas
$ws.Range("A26").Value = "FORMULARIO DE BANCOS IMPLEMENTAR BAJA LOGICA JUNTO CON CUENTAS"
$ws.Range("B26").Value = "no comenzado"

# Nueva tarea (fila 27): seguimiento de OT
$ws.Range("A27").Value = "agregar en seguimiento de OT el estado del OT buscada"
$ws.Range("B27").Value = "no comenzado"

# Actualizar la posicion de scroll/seleccion de la vista de la hoja
$excel.Goto($ws.Range("A10"), $true)
$ws.Range("C28").Select()
